$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Dades_Període")
$ws2 = $wb.Worksheets.Item("Estudi_Capçaleres")

# Update DATA_EXTRACCIO (column H) timestamps for all data rows (2-226) on sheet "Dades_Període"
$ws1.Cells.Item(2, 8).Value = "2026-02-19 07:28:43"
$ws1.Cells.Item(3, 8).Value = "2026-02-19 07:28:44"
$ws1.Cells.Item(4, 8).Value = "2026-02-19 07:28:44"
$ws1.Cells.Item(5, 8).Value = "2026-02-19 07:28:44"
$ws1.Cells.Item(6, 8).Value = "2026-02-19 07:28:44"
$ws1.Cells.Item(7, 8).Value = "2026-02-19 07:28:47"
$ws1.Cells.Item(8, 8).Value = "2026-02-19 07:28:49"
$ws1.Cells.Item(9, 8).Value = "2026-02-19 07:28:49"
$ws1.Cells.Item(10, 8).Value = "2026-02-19 07:28:49"
$ws1.Cells.Item(11, 8).Value = "2026-02-19 07:28:49"
$ws1.Cells.Item(12, 8).Value = "2026-02-19 07:28:52"
$ws1.Cells.Item(13, 8).Value = "2026-02-19 07:28:53"
$ws1.Cells.Item(14, 8).Value = "2026-02-19 07:28:53"
$ws1.Cells.Item(15, 8).Value = "2026-02-19 07:28:53"
$ws1.Cells.Item(16, 8).Value = "2026-02-19 07:28:53"
$ws1.Cells.Item(17, 8).Value = "2026-02-19 07:28:56"
$ws1.Cells.Item(18, 8).Value = "2026-02-19 07:28:58"
$ws1.Cells.Item(19, 8).Value = "2026-02-19 07:28:58"
$ws1.Cells.Item(20, 8).Value = "2026-02-19 07:28:58"
$ws1.Cells.Item(21, 8).Value = "2026-02-19 07:28:58"
$ws1.Cells.Item(22, 8).Value = "2026-02-19 07:29:00"
$ws1.Cells.Item(23, 8).Value = "2026-02-19 07:29:02"
$ws1.Cells.Item(24, 8).Value = "2026-02-19 07:29:02"
$ws1.Cells.Item(25, 8).Value = "2026-02-19 07:29:02"
$ws1.Cells.Item(26, 8).Value = "2026-02-19 07:29:02"
$ws1.Cells.Item(27, 8).Value = "2026-02-19 07:29:05"
$ws1.Cells.Item(28, 8).Value = "2026-02-19 07:29:06"
$ws1.Cells.Item(29, 8).Value = "2026-02-19 07:29:06"
$ws1.Cells.Item(30, 8).Value = "2026-02-19 07:29:06"
$ws1.Cells.Item(31, 8).Value = "2026-02-19 07:29:06"
$ws1.Cells.Item(32, 8).Value = "2026-02-19 07:29:09"
$ws1.Cells.Item(33, 8).Value = "2026-02-19 07:29:11"
$ws1.Cells.Item(34, 8).Value = "2026-02-19 07:29:11"
$ws1.Cells.Item(35, 8).Value = "2026-02-19 07:29:11"
$ws1.Cells.Item(36, 8).Value = "2026-02-19 07:29:11"
$ws1.Cells.Item(37, 8).Value = "2026-02-19 07:29:14"
$ws1.Cells.Item(38, 8).Value = "2026-02-19 07:29:15"
$ws1.Cells.Item(39, 8).Value = "2026-02-19 07:29:15"
$ws1.Cells.Item(40, 8).Value = "2026-02-19 07:29:15"
$ws1.Cells.Item(41, 8).Value = "2026-02-19 07:29:15"
$ws1.Cells.Item(42, 8).Value = "2026-02-19 07:29:18"
$ws1.Cells.Item(43, 8).Value = "2026-02-19 07:29:19"
$ws1.Cells.Item(44, 8).Value = "2026-02-19 07:29:19"
$ws1.Cells.Item(45, 8).Value = "2026-02-19 07:29:19"
$ws1.Cells.Item(46, 8).Value = "2026-02-19 07:29:19"
$ws1.Cells.Item(47, 8).Value = "2026-02-19 07:29:22"
$ws1.Cells.Item(48, 8).Value = "2026-02-19 07:29:24"
$ws1.Cells.Item(49, 8).Value = "2026-02-19 07:29:24"
$ws1.Cells.Item(50, 8).Value = "2026-02-19 07:29:24"
$ws1.Cells.Item(51, 8).Value = "2026-02-19 07:29:24"
$ws1.Cells.Item(52, 8).Value = "2026-02-19 07:29:26"
$ws1.Cells.Item(53, 8).Value = "2026-02-19 07:29:28"
$ws1.Cells.Item(54, 8).Value = "2026-02-19 07:29:28"
$ws1.Cells.Item(55, 8).Value = "2026-02-19 07:29:28"
$ws1.Cells.Item(56, 8).Value = "2026-02-19 07:29:28"
$ws1.Cells.Item(57, 8).Value = "2026-02-19 07:29:30"
$ws1.Cells.Item(58, 8).Value = "2026-02-19 07:29:32"
$ws1.Cells.Item(59, 8).Value = "2026-02-19 07:29:32"
$ws1.Cells.Item(60, 8).Value = "2026-02-19 07:29:32"
$ws1.Cells.Item(61, 8).Value = "2026-02-19 07:29:32"
$ws1.Cells.Item(62, 8).Value = "2026-02-19 07:29:35"
$ws1.Cells.Item(63, 8).Value = "2026-02-19 07:29:36"
$ws1.Cells.Item(64, 8).Value = "2026-02-19 07:29:36"
$ws1.Cells.Item(65, 8).Value = "2026-02-19 07:29:36"
$ws1.Cells.Item(66, 8).Value = "2026-02-19 07:29:36"
$ws1.Cells.Item(67, 8).Value = "2026-02-19 07:29:39"
$ws1.Cells.Item(68, 8).Value = "2026-02-19 07:29:40"
$ws1.Cells.Item(69, 8).Value = "2026-02-19 07:29:40"
$ws1.Cells.Item(70, 8).Value = "2026-02-19 07:29:40"
$ws1.Cells.Item(71, 8).Value = "2026-02-19 07:29:40"
$ws1.Cells.Item(72, 8).Value = "2026-02-19 07:29:43"
$ws1.Cells.Item(73, 8).Value = "2026-02-19 07:29:45"
$ws1.Cells.Item(74, 8).Value = "2026-02-19 07:29:45"
$ws1.Cells.Item(75, 8).Value = "2026-02-19 07:29:45"
$ws1.Cells.Item(76, 8).Value = "2026-02-19 07:29:45"
$ws1.Cells.Item(77, 8).Value = "2026-02-19 07:29:48"
$ws1.Cells.Item(78, 8).Value = "2026-02-19 07:29:49"
$ws1.Cells.Item(79, 8).Value = "2026-02-19 07:29:49"
$ws1.Cells.Item(80, 8).Value = "2026-02-19 07:29:49"
$ws1.Cells.Item(81, 8).Value = "2026-02-19 07:29:49"
$ws1.Cells.Item(82, 8).Value = "2026-02-19 07:29:52"
$ws1.Cells.Item(83, 8).Value = "2026-02-19 07:29:54"
$ws1.Cells.Item(84, 8).Value = "2026-02-19 07:29:54"
$ws1.Cells.Item(85, 8).Value = "2026-02-19 07:29:54"
$ws1.Cells.Item(86, 8).Value = "2026-02-19 07:29:54"
$ws1.Cells.Item(87, 8).Value = "2026-02-19 07:29:56"
$ws1.Cells.Item(88, 8).Value = "2026-02-19 07:29:58"
$ws1.Cells.Item(89, 8).Value = "2026-02-19 07:29:58"
$ws1.Cells.Item(90, 8).Value = "2026-02-19 07:29:58"
$ws1.Cells.Item(91, 8).Value = "2026-02-19 07:29:58"
$ws1.Cells.Item(92, 8).Value = "2026-02-19 07:30:00"
$ws1.Cells.Item(93, 8).Value = "2026-02-19 07:30:02"
$ws1.Cells.Item(94, 8).Value = "2026-02-19 07:30:02"
$ws1.Cells.Item(95, 8).Value = "2026-02-19 07:30:02"
$ws1.Cells.Item(96, 8).Value = "2026-02-19 07:30:02"
$ws1.Cells.Item(97, 8).Value = "2026-02-19 07:30:05"
$ws1.Cells.Item(98, 8).Value = "2026-02-19 07:30:06"
$ws1.Cells.Item(99, 8).Value = "2026-02-19 07:30:06"
$ws1.Cells.Item(100, 8).Value = "2026-02-19 07:30:06"
$ws1.Cells.Item(101, 8).Value = "2026-02-19 07:30:06"
$ws1.Cells.Item(102, 8).Value = "2026-02-19 07:30:09"
$ws1.Cells.Item(103, 8).Value = "2026-02-19 07:30:10"
$ws1.Cells.Item(104, 8).Value = "2026-02-19 07:30:10"
$ws1.Cells.Item(105, 8).Value = "2026-02-19 07:30:10"
$ws1.Cells.Item(106, 8).Value = "2026-02-19 07:30:10"
$ws1.Cells.Item(107, 8).Value = "2026-02-19 07:30:13"
$ws1.Cells.Item(108, 8).Value = "2026-02-19 07:30:15"
$ws1.Cells.Item(109, 8).Value = "2026-02-19 07:30:15"
$ws1.Cells.Item(110, 8).Value = "2026-02-19 07:30:15"
$ws1.Cells.Item(111, 8).Value = "2026-02-19 07:30:15"
$ws1.Cells.Item(112, 8).Value = "2026-02-19 07:30:17"
$ws1.Cells.Item(113, 8).Value = "2026-02-19 07:30:19"
$ws1.Cells.Item(114, 8).Value = "2026-02-19 07:30:19"
$ws1.Cells.Item(115, 8).Value = "2026-02-19 07:30:19"
$ws1.Cells.Item(116, 8).Value = "2026-02-19 07:30:19"
$ws1.Cells.Item(117, 8).Value = "2026-02-19 07:30:22"
$ws1.Cells.Item(118, 8).Value = "2026-02-19 07:30:23"
$ws1.Cells.Item(119, 8).Value = "2026-02-19 07:30:23"
$ws1.Cells.Item(120, 8).Value = "2026-02-19 07:30:23"
$ws1.Cells.Item(121, 8).Value = "2026-02-19 07:30:23"
$ws1.Cells.Item(122, 8).Value = "2026-02-19 07:30:26"
$ws1.Cells.Item(123, 8).Value = "2026-02-19 07:30:28"
$ws1.Cells.Item(124, 8).Value = "2026-02-19 07:30:28"
$ws1.Cells.Item(125, 8).Value = "2026-02-19 07:30:28"
$ws1.Cells.Item(126, 8).Value = "2026-02-19 07:30:28"
$ws1.Cells.Item(127, 8).Value = "2026-02-19 07:30:31"
$ws1.Cells.Item(128, 8).Value = "2026-02-19 07:30:32"
$ws1.Cells.Item(129, 8).Value = "2026-02-19 07:30:32"
$ws1.Cells.Item(130, 8).Value = "2026-02-19 07:30:32"
$ws1.Cells.Item(131, 8).Value = "2026-02-19 07:30:32"
$ws1.Cells.Item(132, 8).Value = "2026-02-19 07:30:35"
$ws1.Cells.Item(133, 8).Value = "2026-02-19 07:30:37"
$ws1.Cells.Item(134, 8).Value = "2026-02-19 07:30:37"
$ws1.Cells.Item(135, 8).Value = "2026-02-19 07:30:37"
$ws1.Cells.Item(136, 8).Value = "2026-02-19 07:30:37"
$ws1.Cells.Item(137, 8).Value = "2026-02-19 07:30:40"
$ws1.Cells.Item(138, 8).Value = "2026-02-19 07:30:41"
$ws1.Cells.Item(139, 8).Value = "2026-02-19 07:30:41"
$ws1.Cells.Item(140, 8).Value = "2026-02-19 07:30:41"
$ws1.Cells.Item(141, 8).Value = "2026-02-19 07:30:41"
$ws1.Cells.Item(142, 8).Value = "2026-02-19 07:30:44"
$ws1.Cells.Item(143, 8).Value = "2026-02-19 07:30:46"
$ws1.Cells.Item(144, 8).Value = "2026-02-19 07:30:46"
$ws1.Cells.Item(145, 8).Value = "2026-02-19 07:30:46"
$ws1.Cells.Item(146, 8).Value = "2026-02-19 07:30:46"
$ws1.Cells.Item(147, 8).Value = "2026-02-19 07:30:48"
$ws1.Cells.Item(148, 8).Value = "2026-02-19 07:30:50"
$ws1.Cells.Item(149, 8).Value = "2026-02-19 07:30:50"
$ws1.Cells.Item(150, 8).Value = "2026-02-19 07:30:50"
$ws1.Cells.Item(151, 8).Value = "2026-02-19 07:30:50"
$ws1.Cells.Item(152, 8).Value = "2026-02-19 07:30:53"
$ws1.Cells.Item(153, 8).Value = "2026-02-19 07:30:54"
$ws1.Cells.Item(154, 8).Value = "2026-02-19 07:30:54"
$ws1.Cells.Item(155, 8).Value = "2026-02-19 07:30:54"
$ws1.Cells.Item(156, 8).Value = "2026-02-19 07:30:54"
$ws1.Cells.Item(157, 8).Value = "2026-02-19 07:30:57"
$ws1.Cells.Item(158, 8).Value = "2026-02-19 07:30:58"
$ws1.Cells.Item(159, 8).Value = "2026-02-19 07:30:58"
$ws1.Cells.Item(160, 8).Value = "2026-02-19 07:30:58"
$ws1.Cells.Item(161, 8).Value = "2026-02-19 07:30:58"
$ws1.Cells.Item(162, 8).Value = "2026-02-19 07:31:01"
$ws1.Cells.Item(163, 8).Value = "2026-02-19 07:31:03"
$ws1.Cells.Item(164, 8).Value = "2026-02-19 07:31:03"
$ws1.Cells.Item(165, 8).Value = "2026-02-19 07:31:03"
$ws1.Cells.Item(166, 8).Value = "2026-02-19 07:31:03"
$ws1.Cells.Item(167, 8).Value = "2026-02-19 07:31:06"
$ws1.Cells.Item(168, 8).Value = "2026-02-19 07:31:08"
$ws1.Cells.Item(169, 8).Value = "2026-02-19 07:31:08"
$ws1.Cells.Item(170, 8).Value = "2026-02-19 07:31:08"
$ws1.Cells.Item(171, 8).Value = "2026-02-19 07:31:08"
$ws1.Cells.Item(172, 8).Value = "2026-02-19 07:31:10"
$ws1.Cells.Item(173, 8).Value = "2026-02-19 07:31:12"
$ws1.Cells.Item(174, 8).Value = "2026-02-19 07:31:12"
$ws1.Cells.Item(175, 8).Value = "2026-02-19 07:31:12"
$ws1.Cells.Item(176, 8).Value = "2026-02-19 07:31:12"
$ws1.Cells.Item(177, 8).Value = "2026-02-19 07:31:15"
$ws1.Cells.Item(178, 8).Value = "2026-02-19 07:31:16"
$ws1.Cells.Item(179, 8).Value = "2026-02-19 07:31:16"
$ws1.Cells.Item(180, 8).Value = "2026-02-19 07:31:16"
$ws1.Cells.Item(181, 8).Value = "2026-02-19 07:31:16"
$ws1.Cells.Item(182, 8).Value = "2026-02-19 07:31:19"
$ws1.Cells.Item(183, 8).Value = "2026-02-19 07:31:20"
$ws1.Cells.Item(184, 8).Value = "2026-02-19 07:31:20"
$ws1.Cells.Item(185, 8).Value = "2026-02-19 07:31:20"
$ws1.Cells.Item(186, 8).Value = "2026-02-19 07:31:20"
$ws1.Cells.Item(187, 8).Value = "2026-02-19 07:31:23"
$ws1.Cells.Item(188, 8).Value = "2026-02-19 07:31:25"
$ws1.Cells.Item(189, 8).Value = "2026-02-19 07:31:25"
$ws1.Cells.Item(190, 8).Value = "2026-02-19 07:31:25"
$ws1.Cells.Item(191, 8).Value = "2026-02-19 07:31:25"
$ws1.Cells.Item(192, 8).Value = "2026-02-19 07:31:28"
$ws1.Cells.Item(193, 8).Value = "2026-02-19 07:31:29"
$ws1.Cells.Item(194, 8).Value = "2026-02-19 07:31:29"
$ws1.Cells.Item(195, 8).Value = "2026-02-19 07:31:29"
$ws1.Cells.Item(196, 8).Value = "2026-02-19 07:31:29"
$ws1.Cells.Item(197, 8).Value = "2026-02-19 07:31:32"
$ws1.Cells.Item(198, 8).Value = "2026-02-19 07:31:34"
$ws1.Cells.Item(199, 8).Value = "2026-02-19 07:31:34"
$ws1.Cells.Item(200, 8).Value = "2026-02-19 07:31:34"
$ws1.Cells.Item(201, 8).Value = "2026-02-19 07:31:34"
$ws1.Cells.Item(202, 8).Value = "2026-02-19 07:31:36"
$ws1.Cells.Item(203, 8).Value = "2026-02-19 07:31:38"
$ws1.Cells.Item(204, 8).Value = "2026-02-19 07:31:38"
$ws1.Cells.Item(205, 8).Value = "2026-02-19 07:31:38"
$ws1.Cells.Item(206, 8).Value = "2026-02-19 07:31:38"
$ws1.Cells.Item(207, 8).Value = "2026-02-19 07:31:41"
$ws1.Cells.Item(208, 8).Value = "2026-02-19 07:31:42"
$ws1.Cells.Item(209, 8).Value = "2026-02-19 07:31:42"
$ws1.Cells.Item(210, 8).Value = "2026-02-19 07:31:42"
$ws1.Cells.Item(211, 8).Value = "2026-02-19 07:31:42"
$ws1.Cells.Item(212, 8).Value = "2026-02-19 07:31:45"
$ws1.Cells.Item(213, 8).Value = "2026-02-19 07:31:47"
$ws1.Cells.Item(214, 8).Value = "2026-02-19 07:31:47"
$ws1.Cells.Item(215, 8).Value = "2026-02-19 07:31:47"
$ws1.Cells.Item(216, 8).Value = "2026-02-19 07:31:47"
$ws1.Cells.Item(217, 8).Value = "2026-02-19 07:31:50"
$ws1.Cells.Item(218, 8).Value = "2026-02-19 07:31:52"
$ws1.Cells.Item(219, 8).Value = "2026-02-19 07:31:52"
$ws1.Cells.Item(220, 8).Value = "2026-02-19 07:31:52"
$ws1.Cells.Item(221, 8).Value = "2026-02-19 07:31:52"
$ws1.Cells.Item(222, 8).Value = "2026-02-19 07:31:54"
$ws1.Cells.Item(223, 8).Value = "2026-02-19 07:31:56"
$ws1.Cells.Item(224, 8).Value = "2026-02-19 07:31:56"
$ws1.Cells.Item(225, 8).Value = "2026-02-19 07:31:56"
$ws1.Cells.Item(226, 8).Value = "2026-02-19 07:31:56"

# For the 26 station "first rows" (I column 06:30 -> 07:00, and J URL T06:30Z -> T07:00Z)
$ws1.Cells.Item(97, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(97, 10).Value2
$ws1.Cells.Item(97, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(102, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(102, 10).Value2
$ws1.Cells.Item(102, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(107, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(107, 10).Value2
$ws1.Cells.Item(107, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(112, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(112, 10).Value2
$ws1.Cells.Item(112, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(117, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(117, 10).Value2
$ws1.Cells.Item(117, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(122, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(122, 10).Value2
$ws1.Cells.Item(122, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(127, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(127, 10).Value2
$ws1.Cells.Item(127, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(132, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(132, 10).Value2
$ws1.Cells.Item(132, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(137, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(137, 10).Value2
$ws1.Cells.Item(137, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(142, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(142, 10).Value2
$ws1.Cells.Item(142, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(147, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(147, 10).Value2
$ws1.Cells.Item(147, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(152, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(152, 10).Value2
$ws1.Cells.Item(152, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(157, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(157, 10).Value2
$ws1.Cells.Item(157, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(162, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(162, 10).Value2
$ws1.Cells.Item(162, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(167, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(167, 10).Value2
$ws1.Cells.Item(167, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(172, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(172, 10).Value2
$ws1.Cells.Item(172, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(177, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(177, 10).Value2
$ws1.Cells.Item(177, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(182, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(182, 10).Value2
$ws1.Cells.Item(182, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(187, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(187, 10).Value2
$ws1.Cells.Item(187, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(192, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(192, 10).Value2
$ws1.Cells.Item(192, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(197, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(197, 10).Value2
$ws1.Cells.Item(197, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(202, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(202, 10).Value2
$ws1.Cells.Item(202, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(207, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(207, 10).Value2
$ws1.Cells.Item(207, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(212, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(212, 10).Value2
$ws1.Cells.Item(212, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(217, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(217, 10).Value2
$ws1.Cells.Item(217, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")
$ws1.Cells.Item(222, 9).Value = "07:00"
$oldUrl = $ws1.Cells.Item(222, 10).Value2
$ws1.Cells.Item(222, 10).Value = $oldUrl.Replace("T06:30Z", "T07:00Z")

# Update sheet "Estudi_Capçaleres" column F (URL_FONT) for rows 21-46
$oldUrl2 = $ws2.Cells.Item(21, 6).Value2
$ws2.Cells.Item(21, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(22, 6).Value2
$ws2.Cells.Item(22, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(23, 6).Value2
$ws2.Cells.Item(23, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(24, 6).Value2
$ws2.Cells.Item(24, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(25, 6).Value2
$ws2.Cells.Item(25, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(26, 6).Value2
$ws2.Cells.Item(26, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(27, 6).Value2
$ws2.Cells.Item(27, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(28, 6).Value2
$ws2.Cells.Item(28, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(29, 6).Value2
$ws2.Cells.Item(29, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(30, 6).Value2
$ws2.Cells.Item(30, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(31, 6).Value2
$ws2.Cells.Item(31, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(32, 6).Value2
$ws2.Cells.Item(32, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(33, 6).Value2
$ws2.Cells.Item(33, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(34, 6).Value2
$ws2.Cells.Item(34, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(35, 6).Value2
$ws2.Cells.Item(35, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(36, 6).Value2
$ws2.Cells.Item(36, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(37, 6).Value2
$ws2.Cells.Item(37, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(38, 6).Value2
$ws2.Cells.Item(38, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(39, 6).Value2
$ws2.Cells.Item(39, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(40, 6).Value2
$ws2.Cells.Item(40, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(41, 6).Value2
$ws2.Cells.Item(41, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(42, 6).Value2
$ws2.Cells.Item(42, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(43, 6).Value2
$ws2.Cells.Item(43, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(44, 6).Value2
$ws2.Cells.Item(44, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(45, 6).Value2
$ws2.Cells.Item(45, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
$oldUrl2 = $ws2.Cells.Item(46, 6).Value2
$ws2.Cells.Item(46, 6).Value = $oldUrl2.Replace("T06:30Z", "T07:00Z")
